$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-10
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08)
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45207
}
